$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix email addresses (typo: domain was glued before the trailing digit)
$ws.Range("G4").Value = "nemreg1es1@mail.com"
$ws.Range("G5").Value = "nemreg1es2@mail.com"
$ws.Range("G6").Value = "nemreg1es3@mail.com"

# Remove the now-unused/incorrect email entries in rows 7 and 8
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# Update the active cell selection
$ws.Range("G7").Select()
